# Refresh the cryptos table (Price / Volume(1h) columns) with the latest scrape values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.809.55"
$ws.Range("E2").Value = "'  +1.05%  "
$ws.Range("D3").Value = "'1.825.31"
$ws.Range("E3").Value = "'  +1.68%  "
$ws.Range("E4").Value = "'  -0.11%  "
$ws.Range("D5").Value = "'228.83"
$ws.Range("E5").Value = "'  +0.82%  "
$ws.Range("D6").Value = "'0.578"
$ws.Range("E6").Value = "'  +4.12%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "'  -0.19%  "
$ws.Range("D8").Value = "'34.92"
$ws.Range("E8").Value = "'  +7.57%  "
$ws.Range("E9").Value = "'  +1.88%  "
$ws.Range("D10").Value = "'0.0702"
$ws.Range("E10").Value = "'  +1.33%  "
$ws.Range("D11").Value = "'0.0952"
$ws.Range("E11").Value = "'  +0.26%  "
$ws.Range("D12").Value = "'2.088.84"
$ws.Range("E12").Value = "'  +1.73%  "
$ws.Range("D13").Value = "'11.51"
$ws.Range("E13").Value = "'  +3.69%  "
$ws.Range("D14").Value = "'1.853.49"
$ws.Range("E14").Value = "'  +3.16%  "
$ws.Range("D15").Value = "'0.647"
$ws.Range("E15").Value = "'  +2.32%  "
$ws.Range("D16").Value = "'34.778.13"
$ws.Range("E16").Value = "'  +1.10%  "
$ws.Range("E17").Value = "'  +2.62%  "
$ws.Range("D18").Value = "'69.23"
$ws.Range("E18").Value = "'  +1.31%  "
$ws.Range("D19").Value = "'247.69"
$ws.Range("E19").Value = "'  +0.48%  "
$ws.Range("E20").Value = "'  +0.33%  "
$ws.Range("D21").Value = "'11.67"
$ws.Range("E21").Value = "'  +5.31%  "
$ws.Range("E22").Value = "'  -0.15%  "
$ws.Range("D23").Value = "'4.20"
$ws.Range("E23").Value = "'  +0.94%  "
$ws.Range("D24").Value = "'174.51"
$ws.Range("E24").Value = "'  +7.26%  "
$ws.Range("D25").Value = "'2.10"
$ws.Range("E25").Value = "'  +1.53%  "
$ws.Range("D26").Value = "'7.52"
$ws.Range("E26").Value = "'  +3.64%  "
$ws.Range("D27").Value = "'16.92"
$ws.Range("E27").Value = "'  +3.08%  "
$ws.Range("E28").Value = "'  +2.43%  "
$ws.Range("E29").Value = "'  -0.27%  "
$ws.Range("D30").Value = "'4.00"
$ws.Range("E30").Value = "'  +2.56%  "
$ws.Range("D31").Value = "'0.0533"
$ws.Range("E31").Value = "'  +2.03%  "
$ws.Range("E32").Value = "'  +2.57%  "
$ws.Range("E33").Value = "'  +0.78%  "
$ws.Range("D34").Value = "'1.86"
$ws.Range("E34").Value = "'  +1.60%  "
$ws.Range("D35").Value = "'2.63"
$ws.Range("E35").Value = "'  +0.94%  "
$ws.Range("D36").Value = "'1.415.98"
$ws.Range("E36").Value = "'  -1.92%  "
$ws.Range("E37").Value = "'  +2.28%  "
$ws.Range("E38").Value = "'  +2.36%  "
$ws.Range("D39").Value = "'0.0193"
$ws.Range("E39").Value = "'  +0.66%  "
$ws.Range("D40").Value = "'85.40"
$ws.Range("E40").Value = "'  +1.56%  "
$ws.Range("D41").Value = "'2.87"
$ws.Range("E41").Value = "'  +4.55%  "
$ws.Range("E42").Value = "'  +3.07%  "
$ws.Range("D43").Value = "'2.40"
$ws.Range("E43").Value = "'  -0.18%  "
$ws.Range("D44").Value = "'13.80"
$ws.Range("E44").Value = "'  -0.09%  "
$ws.Range("E45").Value = "'  +3.00%  "
$ws.Range("D46").Value = "'0.0519"
$ws.Range("E46").Value = "'  -1.20%  "
$ws.Range("D47").Value = "'6.10"
$ws.Range("E47").Value = "'  +0.16%  "
$ws.Range("D48").Value = "'1.989.08"
$ws.Range("E48").Value = "'  +2.02%  "
$ws.Range("D49").Value = "'105.83"
$ws.Range("E49").Value = "'  +0.12%  "

# Rows 50/51 swapped order: BabyDogeCoin now ranks above PaxDollar, with refreshed data.
$ws.Range("B50").Value = "'BabyDogeCoin"
$ws.Range("C50").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "'0.0₆0131"
$ws.Range("E50").Value = "'  +1.38%  "

$ws.Range("B51").Value = "'PaxDollar"
$ws.Range("C51").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "'  -0.10%  "
